$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=10; Result="Home Win"; Resultado="Fallo";   Profit=-0.7; ROI=-100; Fecha="2025-08-30 18:23:58" },
    @{ Row=11; Result="Away Win"; Resultado="Acierto"; Profit=1.76; ROI=110;  Fecha="2025-08-30 18:23:58" },
    @{ Row=12; Result="Home Win"; Resultado="Acierto"; Profit=2.62; ROI=75;   Fecha="2025-08-30 18:23:58" },
    @{ Row=13; Result="Home Win"; Resultado="Acierto"; Profit=2.85; ROI=57;   Fecha="2025-08-30 18:23:58" },
    @{ Row=14; Result="Away Win"; Resultado="Acierto"; Profit=0.5;  ROI=125;  Fecha="2025-08-30 18:23:58" },
    @{ Row=15; Result="Home Win"; Resultado="Acierto"; Profit=2.24; ROI=44;   Fecha="2025-08-30 18:23:58" },
    @{ Row=16; Result="Home Win"; Resultado="Acierto"; Profit=3;    ROI=60;   Fecha="2025-08-30 18:23:58" },
    @{ Row=17; Result="Draw";     Resultado="Fallo";   Profit=-2.6; ROI=-100; Fecha="2025-08-30 18:23:58" },
    @{ Row=18; Result="Home Win"; Resultado="Fallo";   Profit=-0.6; ROI=-100; Fecha="2025-08-30 18:23:58" },
    @{ Row=19; Result="Away Win"; Resultado="Fallo";   Profit=-4;   ROI=-100; Fecha="2025-08-30 18:23:58" },
    @{ Row=20; Result="Home Win"; Resultado="Acierto"; Profit=2.67; ROI=62;   Fecha="2025-08-30 18:23:58" },
    @{ Row=21; Result="Away Win"; Resultado="Acierto"; Profit=1.98; ROI=110;  Fecha="2025-08-30 18:23:58" },
    @{ Row=22; Result="Draw";     Resultado="Fallo";   Profit=-4.4; ROI=-100; Fecha="2025-08-30 18:23:58" },
    @{ Row=23; Result="Away Win"; Resultado="Fallo";   Profit=-2.4; ROI=-100; Fecha="2025-08-30 18:23:58" },
    @{ Row=24; Result="Draw";     Resultado="Fallo";   Profit=-1.6; ROI=-100; Fecha="2025-08-30 18:23:58" },
    @{ Row=25; Result="Draw";     Resultado="Fallo";   Profit=-3.3; ROI=-100; Fecha="2025-08-30 18:23:58" },
    @{ Row=26; Result="Draw";     Resultado="Fallo";   Profit=-1.9; ROI=-100; Fecha="2025-08-30 18:23:58" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 12).Value = "Completed"
    $ws.Cells.Item($row, 13).Value = $r.Result
    $ws.Cells.Item($row, 14).Value = $r.Resultado
    $ws.Cells.Item($row, 15).Value = $r.Profit
    $ws.Cells.Item($row, 16).Value = $r.ROI
    $ws.Cells.Item($row, 17).Value = $r.Fecha
}
